# Generate Report for Handoff
#
# The localization status for e2e\b.md has moved from "handed back / in
# sync" to "ready for handoff" again because the previously handed-back
# file is stale relative to the newest source. This updates the Overview
# rollup plus each per-locale (zh-cn / de-de) detail row for b.md:
#   - Status -> "Ready for handoff"
#   - Content Duplicate -> False
#   - Latest Handoff File / Datetime -> the new b.*.xlf handoff package
#   - Error Detail -> explains why the old handback is stale
# Column P (Error Detail) is also widened so the new message is legible.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is e2e\b.md ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-06 20:51:25"

# --- zh-cn sheet: row 3 is b.md ---
# (leading "'" forces literal text so "False" isn't auto-coerced to a Boolean,
#  matching the source data's shared-string "True"/"False" text cells; then
#  reset .Style so the cell doesn't retain a "quote prefix" number format)
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-06 20:51:20"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1bfe0f1e29f08b49f014f2f3ec490835b7d3f714/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83540961763741ef6195e07f661f87f0ac4d7d5f/e2e/b.md."

# --- de-de sheet: row 3 is b.md ---
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-06 20:51:25"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1bfe0f1e29f08b49f014f2f3ec490835b7d3f714/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83540961763741ef6195e07f661f87f0ac4d7d5f/e2e/b.md."

# Widen column P (Error Detail) on both locale sheets so the long message fits.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
